# xleash, test: enhance recursive real excel.
# + Improve ABCSheet docs.

$wb = $excel.ActiveWorkbook

# Sheet "2" is the first / active sheet in the workbook (tabSelected).
$ws = $wb.Worksheets.Item("2")

# Add a new row of data: B7 = "No Recurse", C7 = "bar" (reusing existing string).
$ws.Range("B7").Value = "No Recurse"
$ws.Range("C7").Value = "bar"

# Update the selected / active cell on this sheet to the newly added cell.
$ws.Activate()
$ws.Range("B7").Select()
